$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 393.38095
$ws.Range("J9").Value = 523
$ws.Range("L9").Value = 523
$ws.Range("N9").Value = -861
$ws.Range("H19").Value = 3000
$ws.Range("J19").Value = 3000
$ws.Range("L19").Value = 3000
$ws.Range("N19").Value = -3350
$ws.Range("H41").Value = 1866.8334
$ws.Range("I41").Value = 2374.75
$ws.Range("K41").Value = 2374.75
$ws.Range("M41").Value = -1934.75
$ws.Range("H49").Value = 1539
$ws.Range("I49").Value = 1799
$ws.Range("K49").Value = 5397
$ws.Range("M49").Value = -5261
$ws.Range("H74").Value = 10222
$ws.Range("I74").Value = 10151.272
$ws.Range("K74").Value = 10151.272
$ws.Range("M74").Value = -9215.272000000001
$ws.Range("H77").Value = 10222
$ws.Range("I77").Value = 10151.272
$ws.Range("K77").Value = 50756.36
$ws.Range("M77").Value = -46076.36
$ws.Range("H107").Value = 518.6875
$ws.Range("I107").Value = 547.8461
$ws.Range("J107").Value = 392.33334
$ws.Range("K107").Value = 547.8461
$ws.Range("L107").Value = 392.33334
$ws.Range("M107").Value = 1372.1539
$ws.Range("N107").Value = -4232.33334
$ws.Range("H115").Value = 527.3333
$ws.Range("I115").Value = 527.3333
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 1581.9999
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -14.99990000000003
$ws.Range("N115").ClearContents()
$ws.Range("H116").Value = 15000
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H127").Value = 819.5
$ws.Range("I127").Value = 819.5
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 2458.5
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 2501.5
$ws.Range("N127").ClearContents()
$ws.Range("H131").Value = 2601.625
$ws.Range("I131").Value = 1122.3
$ws.Range("K131").Value = 3366.9
$ws.Range("M131").Value = 1673.1
$ws.Range("H138").Value = 4911.7617
$ws.Range("I138").Value = 6028.5
$ws.Range("J138").Value = 4649
$ws.Range("K138").Value = 18085.5
$ws.Range("L138").Value = 13947
$ws.Range("M138").Value = -12945.5
$ws.Range("N138").Value = -24227

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 2454.5
$ws.Range("I4").Value = 3387.4
$ws.Range("J4").Value = 899.6667
$ws.Range("K4").Value = 3387.4
$ws.Range("L4").Value = 899.6667
$ws.Range("M4").Value = -3271.4
$ws.Range("N4").Value = -1131.6667
$ws.Range("H15").Value = 14055.556
$ws.Range("I15").Value = 6500
$ws.Range("J15").Value = 15000
$ws.Range("K15").Value = 6500
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = -6150
$ws.Range("N15").Value = -15700
$ws.Range("H32").Value = 9368.424999999999
$ws.Range("I32").Value = 8086.027
$ws.Range("K32").Value = 8086.027
$ws.Range("M32").Value = -7799.027
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H74").Value = 4602630
$ws.Range("I74").Value = 5752206
$ws.Range("K74").Value = 5752206
$ws.Range("M74").Value = -5751332
$ws.Range("H77").Value = 4602630
$ws.Range("I77").Value = 5752206
$ws.Range("K77").Value = 28761030
$ws.Range("M77").Value = -28756662
$ws.Range("H95").Value = 30861
$ws.Range("J95").Value = 34481.332
$ws.Range("L95").Value = 34481.332
$ws.Range("N95").Value = -39973.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 685.75
$ws.Range("I22").Value = 655.1429000000001
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 655.1429000000001
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -482.1429000000001
$ws.Range("N22").Value = -1246
$ws.Range("H107").Value = 1975.2307
$ws.Range("I107").Value = 1975.2307
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1975.2307
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -55.23070000000007
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1585.75
$ws.Range("I105").Value = 1625
$ws.Range("J105").Value = 1311
$ws.Range("K105").Value = 1625
$ws.Range("L105").Value = 1311
$ws.Range("M105").Value = 122
$ws.Range("N105").Value = -4805
$ws.Range("H132").Value = 9097712
$ws.Range("I132").Value = 9097712
$ws.Range("K132").Value = 27293136
$ws.Range("M132").Value = -27290606
$ws.Range("H138").Value = 180000
$ws.Range("J138").Value = 180000
$ws.Range("L138").Value = 180000
$ws.Range("N138").Value = -190280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2639.9
$ws.Range("I11").Value = 2071.2856
$ws.Range("K11").Value = 6213.8568
$ws.Range("M11").Value = -6073.8568
$ws.Range("H37").Value = 98749.75
$ws.Range("J37").Value = 98749.75
$ws.Range("L37").Value = 296249.25
$ws.Range("N37").Value = -296473.25
$ws.Range("H126").Value = 17799.25
$ws.Range("I126").Value = 8267.5
$ws.Range("K126").Value = 24802.5
$ws.Range("M126").Value = -19862.5
$ws.Range("H129").Value = 2291.7
$ws.Range("I129").Value = 1832.6666
$ws.Range("J129").Value = 2488.4285
$ws.Range("K129").Value = 5497.9998
$ws.Range("L129").Value = 7465.2855
$ws.Range("M129").Value = -497.9997999999996
$ws.Range("N129").Value = -17465.2855
$ws.Range("H137").Value = 3496
$ws.Range("J137").Value = 8998.5
$ws.Range("L137").Value = 26995.5
$ws.Range("N137").Value = -37195.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2404.3
$ws.Range("I80").Value = 1752.8572
$ws.Range("J80").Value = 3924.3333
$ws.Range("K80").Value = 1752.8572
$ws.Range("L80").Value = 3924.3333
$ws.Range("M80").Value = -754.8571999999999
$ws.Range("N80").Value = -5920.3333
$ws.Range("H83").Value = 2404.3
$ws.Range("I83").Value = 1752.8572
$ws.Range("J83").Value = 3924.3333
$ws.Range("K83").Value = 8764.286
$ws.Range("L83").Value = 19621.6665
$ws.Range("M83").Value = -3772.286
$ws.Range("N83").Value = -29605.6665
$ws.Range("H132").Value = 37039604
$ws.Range("I132").Value = 3000.6667
$ws.Range("J132").Value = 111112820
$ws.Range("K132").Value = 9002.000100000001
$ws.Range("L132").Value = 333338460
$ws.Range("M132").Value = -6472.000100000001
$ws.Range("N132").Value = -333343520

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 826.7143
$ws.Range("I16").Value = 797.8333
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 797.8333
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -627.8333
$ws.Range("N16").Value = -1340
$ws.Range("H40").Value = 29417674
$ws.Range("I40").Value = 43483376
$ws.Range("K40").Value = 43483376
$ws.Range("M40").Value = -43483240
$ws.Range("H122").Value = 4750
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H136").Value = 2204.0908
$ws.Range("I136").Value = 2104.5
$ws.Range("J136").Value = 3200
$ws.Range("K136").Value = 6313.5
$ws.Range("L136").Value = 9600
$ws.Range("M136").Value = -3763.5
$ws.Range("N136").Value = -14700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 59102.2
$ws.Range("I45").Value = 45159
$ws.Range("J45").Value = 68397.664
$ws.Range("K45").Value = 45159
$ws.Range("L45").Value = 68397.664
$ws.Range("M45").Value = -44668
$ws.Range("N45").Value = -69379.664
$ws.Range("H122").Value = 2252.5925
$ws.Range("I122").Value = 2146.923
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6440.768999999999
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -3990.768999999999
$ws.Range("N122").Value = -19900
